$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1932
$ws.Range("F3").Value = 33
$ws.Range("F5").Value = 408
$ws.Range("F6").Value = 1819
$ws.Range("F7").Value = 855
$ws.Range("F8").Value = 1275
$ws.Range("F9").Value = 67
$ws.Range("F10").Value = 454
$ws.Range("F12").Value = 2651
$ws.Range("F13").Value = 358
$ws.Range("F14").Value = 866
$ws.Range("F15").Value = 1070
$ws.Range("F17").Value = 16
$ws.Range("F18").Value = 54
$ws.Range("F19").Value = 1549
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 1221
$ws.Range("F22").Value = 168
$ws.Range("F24").Value = 1368
$ws.Range("F25").Value = 47
$ws.Range("F26").Value = 1397
$ws.Range("F27").Value = 945
$ws.Range("F28").Value = 1305
$ws.Range("F29").Value = 199
$ws.Range("F30").Value = 1261
$ws.Range("F31").Value = 420
$ws.Range("F32").Value = 137
$ws.Range("F33").Value = 952
$ws.Range("F34").Value = 21
$ws.Range("F35").Value = 1822
$ws.Range("F36").Value = 457
$ws.Range("F37").Value = 33
$ws.Range("F38").Value = 147
$ws.Range("F40").Value = 2238
$ws.Range("F41").Value = 130
$ws.Range("F42").Value = 878
$ws.Range("F43").Value = 2738
$ws.Range("F44").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 59
$ws.Range("F5").Value = 57
$ws.Range("F6").Value = 89
$ws.Range("F7").Value = 27
$ws.Range("F9").Value = 469
$ws.Range("F11").Value = 17
$ws.Range("F12").Value = 357
$ws.Range("F13").Value = 106221
$ws.Range("F15").Value = 11
$ws.Range("F17").Value = 59
$ws.Range("F18").Value = 59
$ws.Range("F19").Value = 210
$ws.Range("F20").Value = 283
$ws.Range("F22").Value = 269
$ws.Range("F23").Value = 64
$ws.Range("F24").Value = 78
$ws.Range("F25").Value = 66
$ws.Range("F29").Value = 44
$ws.Range("F30").Value = 213
$ws.Range("F32").Value = 0

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 3007
$ws.Range("F6").Value = 4838
$ws.Range("F7").Value = 174
$ws.Range("F9").Value = 653
$ws.Range("F10").Value = 913
$ws.Range("F11").Value = 532
$ws.Range("F12").Value = 584
$ws.Range("F13").Value = 1324
$ws.Range("F14").Value = 375
$ws.Range("F15").Value = 1156

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1932
$ws.Range("F4").Value = 33
$ws.Range("F5").Value = 4838
$ws.Range("F6").Value = 653
$ws.Range("F7").Value = 913
$ws.Range("F8").Value = 532
$ws.Range("F9").Value = 584
$ws.Range("F10").Value = 1324
$ws.Range("F11").Value = 408
$ws.Range("F12").Value = 1819
$ws.Range("F13").Value = 855
$ws.Range("F14").Value = 1275
$ws.Range("F15").Value = 27
$ws.Range("F16").Value = 454
$ws.Range("F17").Value = 1156
$ws.Range("F18").Value = 2651
$ws.Range("F20").Value = 358
$ws.Range("F21").Value = 866
$ws.Range("F22").Value = 1070
$ws.Range("F24").Value = 1549
$ws.Range("F25").Value = 17
$ws.Range("F26").Value = 357
$ws.Range("F27").Value = 1221
$ws.Range("F28").Value = 168
$ws.Range("F30").Value = 1397
$ws.Range("F31").Value = 946
$ws.Range("F32").Value = 1305
$ws.Range("F33").Value = 199
$ws.Range("F34").Value = 11
$ws.Range("F35").Value = 59
$ws.Range("F36").Value = 1261
$ws.Range("F37").Value = 420
$ws.Range("F38").Value = 952
$ws.Range("F39").Value = 64
$ws.Range("F40").Value = 1822
$ws.Range("F42").Value = 33
$ws.Range("F43").Value = 147
$ws.Range("F44").Value = 2238
$ws.Range("F45").Value = 130
$ws.Range("F46").Value = 878
$ws.Range("F47").Value = 2738
